$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 3435.566666666667, 3512, 3399, 0.0288921594619751),
    @(1, 3234.333333333333, 3387, 3029, 0.03049564361572266),
    @(2, 3458.2, 3585, 3238, 0.03043496608734131),
    @(3, 3570.733333333333, 3917, 3312, 0.03012724717458089),
    @(4, 3298, 3555, 2981, 0.03100792566935221),
    @(5, 4081.466666666667, 4227, 3990, 0.03061772187550863),
    @(6, 3107.2, 3344, 2966, 0.03026014963785807),
    @(7, 3790.366666666667, 4021, 3582, 0.03086214860280355),
    @(8, 3149.266666666667, 3377, 2935, 0.03164127667744954),
    @(9, 3122.9, 3531, 3029, 0.02850728829701742)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
